# "created the create edit page" - add the two new fields (Special Item /
# Is Active) that back the new create/edit item form: new header columns D
# and E, plus default values for each existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("D1").Value = "Special Item"
$ws.Range("E1").Value = "Is Active"

# Default values for the existing rows: every current item is "not
# special" and "active" by default.
$ws.Range("D2:D7").Value = $false
$ws.Range("E2:E7").Value = $true

# Leave the selection where the editor last left it.
$ws.Range("K8").Select() | Out-Null
